$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnfsf12"
$ws.Cells.Item(2, 3).Value = "Cd163"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 5.742066
$ws.Cells.Item(2, 8).Value = 17.226198
$ws.Cells.Item(2, 9).Value = 0.2447097919555983
$ws.Cells.Item(2, 10).Value = 0.2447097919555983
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.6022393333333333
$ws.Cells.Item(2, 14).Value = 1.806718
$ws.Cells.Item(2, 15).Value = 0.1247313749130817
$ws.Cells.Item(2, 16).Value = 0.1247313749130817
$ws.Cells.Item(2, 17).Value = 3.458097999796
$ws.Cells.Item(2, 18).Value = 31.122881998164
$ws.Cells.Item(2, 19).Value = 0.03052298880531596
$ws.Cells.Item(2, 20).Value = 0.03052298880531597
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnfsf12"
$ws.Cells.Item(3, 3).Value = "Cd163"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 5.742066
$ws.Cells.Item(3, 8).Value = 17.226198
$ws.Cells.Item(3, 9).Value = 0.2447097919555983
$ws.Cells.Item(3, 10).Value = 0.2447097919555983
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.324840666666667
$ws.Cells.Item(3, 14).Value = 9.974522
$ws.Cells.Item(3, 15).Value = 0.6886165096936998
$ws.Cells.Item(3, 16).Value = 0.6886165096936998
$ws.Cells.Item(3, 17).Value = 19.091454547484
$ws.Cells.Item(3, 18).Value = 171.823090927356
$ws.Cells.Item(3, 19).Value = 0.1685112028243355
$ws.Cells.Item(3, 20).Value = 0.1685112028243355
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnfsf12"
$ws.Cells.Item(4, 3).Value = "Cd163"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 5.742066
$ws.Cells.Item(4, 8).Value = 17.226198
$ws.Cells.Item(4, 9).Value = 0.2447097919555983
$ws.Cells.Item(4, 10).Value = 0.2447097919555983
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.9012106666666666
$ws.Cells.Item(4, 14).Value = 2.703632
$ws.Cells.Item(4, 15).Value = 0.1866521153932185
$ws.Cells.Item(4, 16).Value = 0.1866521153932185
$ws.Cells.Item(4, 17).Value = 5.174811127904
$ws.Cells.Item(4, 18).Value = 46.573300151136
$ws.Cells.Item(4, 19).Value = 0.04567560032594683
$ws.Cells.Item(4, 20).Value = 0.04567560032594684
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnfsf12"
$ws.Cells.Item(5, 3).Value = "Cd163"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 6.924657666666666
$ws.Cells.Item(5, 8).Value = 20.773973
$ws.Cells.Item(5, 9).Value = 0.2951083350441702
$ws.Cells.Item(5, 10).Value = 0.2951083350441703
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.6022393333333333
$ws.Cells.Item(5, 14).Value = 1.806718
$ws.Cells.Item(5, 15).Value = 0.1247313749130817
$ws.Cells.Item(5, 16).Value = 0.1247313749130817
$ws.Cells.Item(5, 17).Value = 4.170301216734888
$ws.Cells.Item(5, 18).Value = 37.532710950614
$ws.Cells.Item(5, 19).Value = 0.03680926837836973
$ws.Cells.Item(5, 20).Value = 0.03680926837836975
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnfsf12"
$ws.Cells.Item(6, 3).Value = "Cd163"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 6.924657666666666
$ws.Cells.Item(6, 8).Value = 20.773973
$ws.Cells.Item(6, 9).Value = 0.2951083350441702
$ws.Cells.Item(6, 10).Value = 0.2951083350441703
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.324840666666667
$ws.Cells.Item(6, 14).Value = 9.974522
$ws.Cells.Item(6, 15).Value = 0.6886165096936998
$ws.Cells.Item(6, 16).Value = 0.6886165096936998
$ws.Cells.Item(6, 17).Value = 23.02338341287844
$ws.Cells.Item(6, 18).Value = 207.210450715906
$ws.Cells.Item(6, 19).Value = 0.2032164716596355
$ws.Cells.Item(6, 20).Value = 0.2032164716596355
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnfsf12"
$ws.Cells.Item(7, 3).Value = "Cd163"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 6.924657666666666
$ws.Cells.Item(7, 8).Value = 20.773973
$ws.Cells.Item(7, 9).Value = 0.2951083350441702
$ws.Cells.Item(7, 10).Value = 0.2951083350441703
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.9012106666666666
$ws.Cells.Item(7, 14).Value = 2.703632
$ws.Cells.Item(7, 15).Value = 0.1866521153932185
$ws.Cells.Item(7, 16).Value = 0.1866521153932185
$ws.Cells.Item(7, 17).Value = 6.24057535221511
$ws.Cells.Item(7, 18).Value = 56.16517816993599
$ws.Cells.Item(7, 19).Value = 0.05508259500616505
$ws.Cells.Item(7, 20).Value = 0.05508259500616507
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Tnfsf12"
$ws.Cells.Item(8, 3).Value = "Cd163"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 8.617968666666666
$ws.Cells.Item(8, 8).Value = 25.853906
$ws.Cells.Item(8, 9).Value = 0.3672722186578602
$ws.Cells.Item(8, 10).Value = 0.3672722186578602
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.6022393333333333
$ws.Cells.Item(8, 14).Value = 1.806718
$ws.Cells.Item(8, 15).Value = 0.1247313749130817
$ws.Cells.Item(8, 16).Value = 0.1247313749130817
$ws.Cells.Item(8, 17).Value = 5.190079704500889
$ws.Cells.Item(8, 18).Value = 46.710717340508
$ws.Cells.Item(8, 19).Value = 0.04581036880057289
$ws.Cells.Item(8, 20).Value = 0.0458103688005729
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Tnfsf12"
$ws.Cells.Item(9, 3).Value = "Cd163"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 8.617968666666666
$ws.Cells.Item(9, 8).Value = 25.853906
$ws.Cells.Item(9, 9).Value = 0.3672722186578602
$ws.Cells.Item(9, 10).Value = 0.3672722186578602
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.324840666666667
$ws.Cells.Item(9, 14).Value = 9.974522
$ws.Cells.Item(9, 15).Value = 0.6886165096936998
$ws.Cells.Item(9, 16).Value = 0.6886165096936998
$ws.Cells.Item(9, 17).Value = 28.65337268699244
$ws.Cells.Item(9, 18).Value = 257.880354182932
$ws.Cells.Item(9, 19).Value = 0.252909713319637
$ws.Cells.Item(9, 20).Value = 0.2529097133196371
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Tnfsf12"
$ws.Cells.Item(10, 3).Value = "Cd163"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 8.617968666666666
$ws.Cells.Item(10, 8).Value = 25.853906
$ws.Cells.Item(10, 9).Value = 0.3672722186578602
$ws.Cells.Item(10, 10).Value = 0.3672722186578602
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.9012106666666666
$ws.Cells.Item(10, 14).Value = 2.703632
$ws.Cells.Item(10, 15).Value = 0.1866521153932185
$ws.Cells.Item(10, 16).Value = 0.1866521153932185
$ws.Cells.Item(10, 17).Value = 7.76660528739911
$ws.Cells.Item(10, 18).Value = 69.89944758659199
$ws.Cells.Item(10, 19).Value = 0.0685521365376503
$ws.Cells.Item(10, 20).Value = 0.06855213653765031
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Tnfsf12"
$ws.Cells.Item(11, 3).Value = "Cd163"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.180106333333333
$ws.Cells.Item(11, 8).Value = 6.540318999999999
$ws.Cells.Item(11, 9).Value = 0.09290965434237122
$ws.Cells.Item(11, 10).Value = 0.09290965434237124
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.6022393333333333
$ws.Cells.Item(11, 14).Value = 1.806718
$ws.Cells.Item(11, 15).Value = 0.1247313749130817
$ws.Cells.Item(11, 16).Value = 0.1247313749130817
$ws.Cells.Item(11, 17).Value = 1.312945784782444
$ws.Cells.Item(11, 18).Value = 11.816512063042
$ws.Cells.Item(11, 19).Value = 0.01158874892882314
$ws.Cells.Item(11, 20).Value = 0.01158874892882314
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Tnfsf12"
$ws.Cells.Item(12, 3).Value = "Cd163"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 2.180106333333333
$ws.Cells.Item(12, 8).Value = 6.540318999999999
$ws.Cells.Item(12, 9).Value = 0.09290965434237122
$ws.Cells.Item(12, 10).Value = 0.09290965434237124
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.324840666666667
$ws.Cells.Item(12, 14).Value = 9.974522
$ws.Cells.Item(12, 15).Value = 0.6886165096936998
$ws.Cells.Item(12, 16).Value = 0.6886165096936998
$ws.Cells.Item(12, 17).Value = 7.248506194724222
$ws.Cells.Item(12, 18).Value = 65.236555752518
$ws.Cells.Item(12, 19).Value = 0.06397912189009178
$ws.Cells.Item(12, 20).Value = 0.06397912189009178
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Tnfsf12"
$ws.Cells.Item(13, 3).Value = "Cd163"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 2.180106333333333
$ws.Cells.Item(13, 8).Value = 6.540318999999999
$ws.Cells.Item(13, 9).Value = 0.09290965434237122
$ws.Cells.Item(13, 10).Value = 0.09290965434237124
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.9012106666666666
$ws.Cells.Item(13, 14).Value = 2.703632
$ws.Cells.Item(13, 15).Value = 0.1866521153932185
$ws.Cells.Item(13, 16).Value = 0.1866521153932185
$ws.Cells.Item(13, 17).Value = 1.964735082067555
$ws.Cells.Item(13, 18).Value = 17.682615738608
$ws.Cells.Item(13, 19).Value = 0.01734178352345632
$ws.Cells.Item(13, 20).Value = 0.01734178352345632
